# [Kadastro App] Yeni kayit eklendi: 2966
#
# Appends one new record (Kayit No 2966) as the next row after the last
# existing row of data on both the master "Kayitlar" sheet and the
# district-filtered "Erdemli" sheet (new row 34 on each, following the
# existing row 33).
#
# All of the existing cells in these tables are stored as literal TEXT
# (even the ones that look like numbers or dates), so every new cell is
# written the same way: as a text value, not a number/date, to keep the
# column's data type consistent with the rest of the sheet.

$wb = $excel.ActiveWorkbook

$newRecord = @{
    KayitNo    = "2966"
    Tarih      = "2025-09-10"
    Birim      = "Erdemli"
    Parsel     = "1"
    Is         = "ÇAP"
    Personeller = "AYHAN KARADAYI (K.Teknisyeni)"
}

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Next free row is right after the current last row of data.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Force each cell to be stored as text (quote-prefixed literal) so the
    # numeric-looking / date-looking values ("2966", "2025-09-10", "1")
    # don't get auto-converted into a number or a date serial, matching
    # how every other row in this table is stored.
    $ws.Cells.Item($newRow, 1).Value = "'" + $newRecord.KayitNo
    $ws.Cells.Item($newRow, 2).Value = "'" + $newRecord.Tarih
    $ws.Cells.Item($newRow, 3).Value = $newRecord.Birim
    $ws.Cells.Item($newRow, 4).Value = "'" + $newRecord.Parsel
    $ws.Cells.Item($newRow, 5).Value = $newRecord.Is
    $ws.Cells.Item($newRow, 6).Value = $newRecord.Personeller
}
